# Generate Report for Archive
#
# Two of the four tracked markdown files (4f25c581-...md and
# 5156710d-...md) have moved from "Ready for handoff" to
# "In Translation" for both target locales (zh-cn, de-de). Update the
# Status column on the per-locale sheets and the corresponding summary
# columns on the Overview sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows for the two files ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus   # 4f25c581-...md / zh-cn
$overview.Range("F3").Value = $newStatus   # 4f25c581-...md / de-de
$overview.Range("E4").Value = $newStatus   # 5156710d-...md / zh-cn
$overview.Range("F4").Value = $newStatus   # 5156710d-...md / de-de

# --- zh-cn sheet: column C (Status) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus       # 4f25c581-...md
$zhcn.Range("C4").Value = $newStatus       # 5156710d-...md

# --- de-de sheet: column C (Status) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus       # 4f25c581-...md
$dede.Range("C4").Value = $newStatus       # 5156710d-...md
